# Apply changes described in the diff to the active workbook/worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric data rows (no shared strings involved) ---

# Row 5
$ws.Range("A5").Value = 20250312
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5

# Row 7
$ws.Range("A7").Value = 20250328
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 6

# Row 8
$ws.Range("G8").Value = 17

# Row 9
$ws.Range("A9").Value = 20250402
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 5

# Row 11
$ws.Range("A11").Value = 20250501
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 5

# Row 12 (numeric part)
$ws.Range("E12").Value = 28

# Row 13
$ws.Range("A13").Value = 20250517
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 4

# Row 15
$ws.Range("A15").Value = 20250519
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5

# Row 17
$ws.Range("A17").Value = 20250520
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4

# Row 18 (numeric part)
$ws.Range("F18").Value = 23

# --- Text / shared-string cells -------------------------------------------
# Entered in the exact order the original author typed them so the
# resulting shared-string table ordering matches the source workbook.

$ws.Range("C10").Value = "15,19,22,27"
$ws.Range("C12").Value = "5,12,19,25"
$ws.Range("B18").Value = "18,23,27"
$ws.Range("G18").Value = "6,12,15,16,20,22,27,28"
$ws.Range("C18").Value = "16,23"
$ws.Range("E18").Value = "9,14,21"
$ws.Range("B14").Value = "7,11,23"
$ws.Range("F14").Value = "10,13,16,20"
$ws.Range("D14").Value = "1,3,9,10,11,12,22,26"
$ws.Range("G14").Value = "4,5,11,14,18,19,25,26"
$ws.Range("E14").Value = "6,11,14,16,18,22,26"
$ws.Range("C14").Value = "2,7,9,10,11,18,22,23,24"
$ws.Range("D16").Value = "9,10,11,12,22,23,28"
$ws.Range("F16").Value = "15,23,24,25,26"
$ws.Range("B16").Value = "3,4,5,9,12,13,14,19,22,23,24,25"
$ws.Range("E16").Value = " 9,10,15,21,22,23,24"
$ws.Range("G16").Value = "1,2,3,8,13,15,18,19,20,26,27,28"
$ws.Range("C16").Value = "1,2,3,9,10,12,17,21,25,26"

# --- Column widths ----------------------------------------------------------
# (values chosen so the engine's internal rounding reproduces the target
# stored widths from the diff as closely as possible)
$ws.Columns.Item(2).ColumnWidth = 26.666666666666668
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332
$ws.Columns.Item(7).ColumnWidth = 22.5

# --- Selection ---------------------------------------------------------------
$ws.Range("C16").Select() | Out-Null
